$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue 'D2' '60.698.06'
Set-TextValue 'E2' '  -0.43%  '

Set-TextValue 'D3' '2.368.82'
Set-TextValue 'E3' '  -3.33%  '

Set-TextValue 'E4' '  +0.15%  '

Set-TextValue 'D5' '542.87'
Set-TextValue 'E5' '  -1.09%  '

Set-TextValue 'D6' '141.04'
Set-TextValue 'E6' '  -2.78%  '

Set-TextValue 'E7' '  +0.13%  '

Set-TextValue 'D8' '0.540'
Set-TextValue 'E8' '  -9.71%  '

Set-TextValue 'D9' '2.368.39'
Set-TextValue 'E9' '  -3.36%  '

Set-TextValue 'E10' '  -1.89%  '

Set-TextValue 'E11' '  +0.55%  '

Set-TextValue 'E12' '  -0.77%  '

Set-TextValue 'D13' '0.344'
Set-TextValue 'E13' '  -2.30%  '

Set-TextValue 'D14' '25.47'
Set-TextValue 'E14' '  -1.87%  '

Set-TextValue 'D15' '2.798.60'
Set-TextValue 'E15' '  -2.98%  '

Set-TextValue 'E16' '  -0.78%  '

Set-TextValue 'D17' '60.470.21'
Set-TextValue 'E17' '  -0.61%  '

Set-TextValue 'D18' '2.371.18'
Set-TextValue 'E18' '  -3.08%  '

Set-TextValue 'D19' '10.66'
Set-TextValue 'E19' '  -3.43%  '

Set-TextValue 'D20' '4.10'
Set-TextValue 'E20' '  -1.44%  '

Set-TextValue 'D21' '316.43'
Set-TextValue 'E21' '  -0.48%  '

Set-TextValue 'D22' '6.69'
Set-TextValue 'E22' '  -2.73%  '

Set-TextValue 'E23' '  -0.13%  '

Set-TextValue 'E24' '  +3.01%  '

Set-TextValue 'D25' '62.92'
Set-TextValue 'E25' '  -0.75%  '

Set-TextValue 'D26' '1.00'
Set-TextValue 'E26' '  +0.22%  '

Set-TextValue 'B27' 'Aptos'
Set-TextValue 'C27' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D27' '7.82'
Set-TextValue 'E27' '  +3.39%  '

Set-TextValue 'B28' 'WrappedeETH'
Set-TextValue 'C28' 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue 'D28' '2.488.00'
Set-TextValue 'E28' '  -3.31%  '

Set-TextValue 'B29' 'PEPE'
Set-TextValue 'C29' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D29' '0.0₃0928'
Set-TextValue 'E29' '  -4.42%  '

Set-TextValue 'D30' '519.94'
Set-TextValue 'E30' '  -3.68%  '

Set-TextValue 'D31' '1.42'
Set-TextValue 'E31' '  -4.73%  '

Set-TextValue 'D32' '8.00'
Set-TextValue 'E32' '  -3.92%  '

Set-TextValue 'E33' '  -3.42%  '

Set-TextValue 'D34' '1.83'
Set-TextValue 'E34' '  -3.09%  '

Set-TextValue 'E35' '  -0.05%  '

Set-TextValue 'E36' '  +0.19%  '

Set-TextValue 'D37' '4.64'
Set-TextValue 'E37' '  -4.00%  '

Set-TextValue 'D38' '5.44'
Set-TextValue 'E38' '  -6.40%  '

Set-TextValue 'E39' '  -0.63%  '

Set-TextValue 'D40' '18.00'
Set-TextValue 'E40' '  -2.30%  '

Set-TextValue 'E41' '  +0.33%  '

Set-TextValue 'E42' '  +0.03%  '

Set-TextValue 'D43' '136.91'
Set-TextValue 'E43' '  -5.37%  '

Set-TextValue 'D44' '40.21'
Set-TextValue 'E44' '  +1.12%  '

Set-TextValue 'D45' '2.21'
Set-TextValue 'E45' '  -3.24%  '

Set-TextValue 'D46' '139.50'
Set-TextValue 'E46' '  -4.21%  '

Set-TextValue 'E47' '  +0.05%  '

Set-TextValue 'D48' '20.25'
Set-TextValue 'E48' '  -2.54%  '

Set-TextValue 'D49' '0.0517'
Set-TextValue 'E49' '  -2.29%  '

Set-TextValue 'D50' '0.574'
Set-TextValue 'E50' '  -1.48%  '

Set-TextValue 'D51' '0.0910'
Set-TextValue 'E51' '  -3.07%  '

